$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "58.166.31"
Set-TextValue "E2" "  -1.73%  "
Set-TextValue "D3" "2.471.31"
Set-TextValue "E3" "  -2.08%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "519.88"
Set-TextValue "E5" "  -3.16%  "
Set-TextValue "D6" "132.47"
Set-TextValue "E6" "  -3.96%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.04%  "
Set-TextValue "E8" "  -1.73%  "
Set-TextValue "D9" "0.0992"
Set-TextValue "E9" "  -2.25%  "
Set-TextValue "E10" "  -0.53%  "
Set-TextValue "D11" "5.37"
Set-TextValue "E11" "  +0.23%  "
Set-TextValue "D12" "0.342"
Set-TextValue "E12" "  -2.07%  "
Set-TextValue "D13" "2.910.21"
Set-TextValue "E13" "  -2.13%  "
Set-TextValue "D14" "58.107.65"
Set-TextValue "E14" "  -1.73%  "
Set-TextValue "D15" "22.11"
Set-TextValue "E15" "  -4.44%  "
Set-TextValue "E16" "  -2.30%  "
Set-TextValue "D17" "2.471.00"
Set-TextValue "E17" "  -2.17%  "
Set-TextValue "D18" "10.85"
Set-TextValue "E18" "  -2.49%  "
Set-TextValue "D19" "4.18"
Set-TextValue "E19" "  -2.78%  "
Set-TextValue "D20" "319.78"
Set-TextValue "E20" "  -1.78%  "
Set-TextValue "E21" "  -0.07%  "
Set-TextValue "D22" "5.75"
Set-TextValue "E22" "  -4.04%  "
Set-TextValue "D23" "64.53"
Set-TextValue "E23" "  -2.47%  "
Set-TextValue "D24" "0.408"
Set-TextValue "E24" "  -3.76%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  -0.17%  "
Set-TextValue "E26" "  -3.53%  "
Set-TextValue "D27" "7.41"
Set-TextValue "E28" "  -3.12%  "
Set-TextValue "D29" "6.39"
Set-TextValue "E29" "  -5.26%  "
Set-TextValue "E30" "  -4.84%  "
Set-TextValue "D31" "165.41"
Set-TextValue "E31" "  +1.74%  "
Set-TextValue "E32" "  -3.99%  "
Set-TextValue "E33" "  +0.00%  "
Set-TextValue "E34" "  -0.10%  "
Set-TextValue "D35" "18.14"
Set-TextValue "E35" "  -1.96%  "
Set-TextValue "D36" "1.33"
Set-TextValue "E36" "  -9.35%  "
Set-TextValue "D37" "3.99"
Set-TextValue "E37" "  -3.40%  "
Set-TextValue "E38" "  -3.78%  "
Set-TextValue "D39" "0.796"
Set-TextValue "E39" "  -2.64%  "
Set-TextValue "D40" "276.68"
Set-TextValue "E40" "  -3.74%  "
Set-TextValue "E41" "  -4.84%  "
Set-TextValue "D42" "5.04"
Set-TextValue "E42" "  -3.47%  "
Set-TextValue "D43" "0.595"
Set-TextValue "E43" "  -2.71%  "
Set-TextValue "D44" "126.20"
Set-TextValue "E44" "  -4.75%  "
Set-TextValue "D45" "0.0908"
Set-TextValue "E45" "  -2.62%  "
Set-TextValue "E46" "  -3.67%  "
Set-TextValue "E47" "  -3.50%  "
Set-TextValue "D48" "17.15"
Set-TextValue "E48" "  -1.44%  "
Set-TextValue "D49" "1.737.93"
Set-TextValue "E49" "  -1.45%  "
Set-TextValue "E50" "  -1.57%  "
Set-TextValue "E51" "  -1.86%  "
